# Update the dSF (column F) values for specific rows to reflect the
# repulled data / recalculated mean (per commit message: "repull data,
# push all data, mean calculation").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F10").Value = -3
$ws.Range("F22").Value = -1
$ws.Range("F24").Value = -1
$ws.Range("F28").Value = -1
$ws.Range("F29").Value = 1
$ws.Range("F30").Value = 4
$ws.Range("F32").Value = 2
